# Weekly update: insert two new rows of fresh price data at the top of the
# data block (rows 123-124), pushing the existing history down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 123 (existing rows 123:240 shift to 125:242)
$ws.Rows("123:124").Insert()

# Row 123: new "Primera" quality entry
$ws.Cells.Item(123, 1).Value = 1
$ws.Cells.Item(123, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(123, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(123, 4).Value = 44554
$ws.Cells.Item(123, 5).Value = 15
$ws.Cells.Item(123, 6).Value = 100112032
$ws.Cells.Item(123, 7).Value = "Zapallo italiano"
$ws.Cells.Item(123, 8).Value = "Huracán"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 120
$ws.Cells.Item(123, 11).Value = 4000
$ws.Cells.Item(123, 12).Value = 4500
$ws.Cells.Item(123, 13).Value = 4250
$ws.Cells.Item(123, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(123, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(123, 16).Value = 61
$ws.Cells.Item(123, 17).Value = 70
$ws.Cells.Item(123, 18).Value = "Hortaliza"

# Row 124: new "Segunda" quality entry
$ws.Cells.Item(124, 1).Value = 1
$ws.Cells.Item(124, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(124, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(124, 4).Value = 44554
$ws.Cells.Item(124, 5).Value = 15
$ws.Cells.Item(124, 6).Value = 100112032
$ws.Cells.Item(124, 7).Value = "Zapallo italiano"
$ws.Cells.Item(124, 8).Value = "Huracán"
$ws.Cells.Item(124, 9).Value = "Segunda"
$ws.Cells.Item(124, 10).Value = 120
$ws.Cells.Item(124, 11).Value = 3000
$ws.Cells.Item(124, 12).Value = 3500
$ws.Cells.Item(124, 13).Value = 3250
$ws.Cells.Item(124, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(124, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(124, 16).Value = 32
$ws.Cells.Item(124, 17).Value = 100
$ws.Cells.Item(124, 18).Value = "Hortaliza"
